$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Cell = "C2"; Value = 1.007233763412072 },
    @{ Cell = "D2"; Value = 1.00991325466653 },
    @{ Cell = "E2"; Value = 1.009822179915485 },
    @{ Cell = "F2"; Value = 1.005418171149548 },
    @{ Cell = "J2"; Value = 1.012506593353203 },
    @{ Cell = "K2"; Value = 1.012785033388608 },
    @{ Cell = "L2"; Value = 1.012694234278073 },
    @{ Cell = "M2"; Value = 1.008303619376537 },
    @{ Cell = "N2"; Value = 1.007964755916351 },
    @{ Cell = "C3"; Value = 1.009287095432225 },
    @{ Cell = "D3"; Value = 1.011764087100809 },
    @{ Cell = "E3"; Value = 1.011604225035833 },
    @{ Cell = "F3"; Value = 1.008140331328462 },
    @{ Cell = "J3"; Value = 1.014185389888848 },
    @{ Cell = "K3"; Value = 1.014436890927127 },
    @{ Cell = "L3"; Value = 1.014277476572159 },
    @{ Cell = "M3"; Value = 1.010823321684792 },
    @{ Cell = "N3"; Value = 1.008557101822972 },
    @{ Cell = "C4"; Value = 1.010609914077341 },
    @{ Cell = "D4"; Value = 1.012956425105848 },
    @{ Cell = "E4"; Value = 1.012752198477129 },
    @{ Cell = "F4"; Value = 1.009894655252363 },
    @{ Cell = "J4"; Value = 1.015265863463752 },
    @{ Cell = "K4"; Value = 1.015500096321514 },
    @{ Cell = "L4"; Value = 1.015296411973175 },
    @{ Cell = "M4"; Value = 1.012446481218014 },
    @{ Cell = "N4"; Value = 1.008937390812613 },
    @{ Cell = "C5"; Value = 1.011164667208386 },
    @{ Cell = "D5"; Value = 1.013456451286343 },
    @{ Cell = "E5"; Value = 1.013233607139252 },
    @{ Cell = "F5"; Value = 1.010630526950307 },
    @{ Cell = "J5"; Value = 1.015718730223423 },
    @{ Cell = "K5"; Value = 1.015945741712944 },
    @{ Cell = "L5"; Value = 1.015723475712681 },
    @{ Cell = "M5"; Value = 1.013127169341355 },
    @{ Cell = "N5"; Value = 1.009096556300793 },
    @{ Cell = "C6"; Value = 1.011257733956608 },
    @{ Cell = "D6"; Value = 1.01354033646001 },
    @{ Cell = "E6"; Value = 1.013314368262514 },
    @{ Cell = "F6"; Value = 1.01075398803822 },
    @{ Cell = "J6"; Value = 1.015794689376723 },
    @{ Cell = "K6"; Value = 1.016020490579834 },
    @{ Cell = "L6"; Value = 1.015795106326709 },
    @{ Cell = "M6"; Value = 1.013241362293346 },
    @{ Cell = "N6"; Value = 1.009123239652299 },
    @{ Cell = "C7"; Value = 1.01061733201275 },
    @{ Cell = "D7"; Value = 1.012963111283718 },
    @{ Cell = "E7"; Value = 1.012758635757744 },
    @{ Cell = "F7"; Value = 1.009904494403763 },
    @{ Cell = "J7"; Value = 1.015271920013607 },
    @{ Cell = "K7"; Value = 1.015506056230265 },
    @{ Cell = "L7"; Value = 1.015302123479047 },
    @{ Cell = "M7"; Value = 1.01245558317387 },
    @{ Cell = "N7"; Value = 1.008939520358412 },
    @{ Cell = "C8"; Value = 1.007928923536867 },
    @{ Cell = "D8"; Value = 1.010539862416146 },
    @{ Cell = "E8"; Value = 1.010425509666389 },
    @{ Cell = "F8"; Value = 1.006339639084469 },
    @{ Cell = "J8"; Value = 1.013075173167773 },
    @{ Cell = "K8"; Value = 1.01334447472508 },
    @{ Cell = "L8"; Value = 1.013230459314931 },
    @{ Cell = "M8"; Value = 1.009156696789252 },
    @{ Cell = "N8"; Value = 1.008165568793857 },
    @{ Cell = "C9"; Value = 1.003145426914801 },
    @{ Cell = "D9"; Value = 1.006228032427489 },
    @{ Cell = "E9"; Value = 1.00627366806349 },
    @{ Cell = "F9"; Value = 1.000001172251525 },
    @{ Cell = "J9"; Value = 1.009158333178677 },
    @{ Cell = "K9"; Value = 1.00949089024572 },
    @{ Cell = "L9"; Value = 1.009536367171991 },
    @{ Cell = "M9"; Value = 1.0032858352567 },
    @{ Cell = "N9"; Value = 1.006778356089765 },
    @{ Cell = "C10"; Value = 0.9999231654433739 },
    @{ Cell = "D10"; Value = 1.003323481638484 },
    @{ Cell = "E10"; Value = 1.00347666406088 },
    @{ Cell = "F10"; Value = 0.9957340015561779 },
    @{ Cell = "J10"; Value = 1.006514390314371 },
    @{ Cell = "K10"; Value = 1.00689004150809 },
    @{ Cell = "L10"; Value = 1.007042636510944 },
    @{ Cell = "M10"; Value = 0.9993299166240083 },
    @{ Cell = "N10"; Value = 1.005837174268722 },
    @{ Cell = "C11"; Value = 0.9985195009428358 },
    @{ Cell = "D11"; Value = 1.002058230093307 },
    @{ Cell = "E11"; Value = 1.002258213239049 },
    @{ Cell = "F11"; Value = 0.9938756573484836 },
    @{ Cell = "J11"; Value = 1.005361353862241 },
    @{ Cell = "K11"; Value = 1.005755897172606 },
    @{ Cell = "L11"; Value = 1.005955082248306 },
    @{ Cell = "M11"; Value = 0.9976062853082815 },
    @{ Cell = "N11"; Value = 1.005425605319509 },
    @{ Cell = "C12"; Value = 0.9979968113523655 },
    @{ Cell = "D12"; Value = 1.001587085480181 },
    @{ Cell = "E12"; Value = 1.001804489229527 },
    @{ Cell = "F12"; Value = 0.9931837220361913 },
    @{ Cell = "J12"; Value = 1.004931797492528 },
    @{ Cell = "K12"; Value = 1.005333394152456 },
    @{ Cell = "L12"; Value = 1.005549917675708 },
    @{ Cell = "M12"; Value = 0.9969643836886336 },
    @{ Cell = "N12"; Value = 1.005272112154014 },
    @{ Cell = "C13"; Value = 0.9981089898545343 },
    @{ Cell = "D13"; Value = 1.00168820133247 },
    @{ Cell = "E13"; Value = 1.001901866611172 },
    @{ Cell = "F13"; Value = 0.9933322208667903 },
    @{ Cell = "J13"; Value = 1.005023996734403 },
    @{ Cell = "K13"; Value = 1.005424078780855 },
    @{ Cell = "L13"; Value = 1.005636881635681 },
    @{ Cell = "M13"; Value = 0.9971021502599338 },
    @{ Cell = "N13"; Value = 1.00530506513454 },
    @{ Cell = "C14"; Value = 0.9984763221568793 },
    @{ Cell = "D14"; Value = 1.002019309300046 },
    @{ Cell = "E14"; Value = 1.002220731678253 },
    @{ Cell = "F14"; Value = 0.9938184960838866 },
    @{ Cell = "J14"; Value = 1.005325872665321 },
    @{ Cell = "K14"; Value = 1.005720998267354 },
    @{ Cell = "L14"; Value = 1.005921615865152 },
    @{ Cell = "M14"; Value = 0.9975532599130635 },
    @{ Cell = "N14"; Value = 1.005412930209645 },
    @{ Cell = "C15"; Value = 0.9987024733256594 },
    @{ Cell = "D15"; Value = 1.002223159157143 },
    @{ Cell = "E15"; Value = 1.00241704320311 },
    @{ Cell = "F15"; Value = 0.9941178838525471 },
    @{ Cell = "J15"; Value = 1.005511699354269 },
    @{ Cell = "K15"; Value = 1.005903775928104 },
    @{ Cell = "L15"; Value = 1.006096890186832 },
    @{ Cell = "M15"; Value = 0.9978309805352875 },
    @{ Cell = "N15"; Value = 1.005479307154862 },
    @{ Cell = "C16"; Value = 1.000016141277716 },
    @{ Cell = "D16"; Value = 1.003407289689546 },
    @{ Cell = "E16"; Value = 1.003557371151442 },
    @{ Cell = "F16"; Value = 0.9958571041249782 },
    @{ Cell = "J16"; Value = 1.006590737850094 },
    @{ Cell = "K16"; Value = 1.006965140220007 },
    @{ Cell = "L16"; Value = 1.007114647702512 },
    @{ Cell = "M16"; Value = 0.9994440778352315 },
    @{ Cell = "N16"; Value = 1.005864402752942 },
    @{ Cell = "C17"; Value = 1.000837890047494 },
    @{ Cell = "D17"; Value = 1.004148012526776 },
    @{ Cell = "E17"; Value = 1.004270681252198 },
    @{ Cell = "F17"; Value = 0.9969451783998144 },
    @{ Cell = "J17"; Value = 1.00726537143314 },
    @{ Cell = "K17"; Value = 1.007628750151266 },
    @{ Cell = "L17"; Value = 1.007750960917954 },
    @{ Cell = "M17"; Value = 1.000453025260632 },
    @{ Cell = "N17"; Value = 1.006104875074253 },
    @{ Cell = "C18"; Value = 1.001316394315355 },
    @{ Cell = "D18"; Value = 1.004579336601009 },
    @{ Cell = "E18"; Value = 1.004686038224455 },
    @{ Cell = "F18"; Value = 0.9975788125116326 },
    @{ Cell = "J18"; Value = 1.007658085910058 },
    @{ Cell = "K18"; Value = 1.008015057034475 },
    @{ Cell = "L18"; Value = 1.008121365968805 },
    @{ Cell = "M18"; Value = 1.001040499793651 },
    @{ Cell = "N18"; Value = 1.006244750435481 },
    @{ Cell = "C19"; Value = 1.001479416177004 },
    @{ Cell = "D19"; Value = 1.00472628477985 },
    @{ Cell = "E19"; Value = 1.004827545758709 },
    @{ Cell = "F19"; Value = 0.9977946941425103 },
    @{ Cell = "J19"; Value = 1.007791858856405 },
    @{ Cell = "K19"; Value = 1.00814664890886 },
    @{ Cell = "L19"; Value = 1.0082475390508 },
    @{ Cell = "M19"; Value = 1.001240641009048 },
    @{ Cell = "N19"; Value = 1.006292378869608 },
    @{ Cell = "C20"; Value = 1.000749808057148 },
    @{ Cell = "D20"; Value = 1.004068615447174 },
    @{ Cell = "E20"; Value = 1.004194222981725 },
    @{ Cell = "F20"; Value = 0.996828544431229 },
    @{ Cell = "J20"; Value = 1.007193071390838 },
    @{ Cell = "K20"; Value = 1.007557630526282 },
    @{ Cell = "L20"; Value = 1.007682767903857 },
    @{ Cell = "M20"; Value = 1.000344881465635 },
    @{ Cell = "N20"; Value = 1.0060791149003 },
    @{ Cell = "C21"; Value = 0.9983681883574549 },
    @{ Cell = "D21"; Value = 1.001921838967306 },
    @{ Cell = "E21"; Value = 1.00212686554567 },
    @{ Cell = "F21"; Value = 0.9936753466134783 },
    @{ Cell = "J21"; Value = 1.005237013008494 },
    @{ Cell = "K21"; Value = 1.005633597157729 },
    @{ Cell = "L21"; Value = 1.005837802090326 },
    @{ Cell = "M21"; Value = 0.997420465916857 },
    @{ Cell = "N21"; Value = 1.005381183798787 },
    @{ Cell = "C22"; Value = 0.9968631906893449 },
    @{ Cell = "D22"; Value = 1.000565264582736 },
    @{ Cell = "E22"; Value = 1.000820438517567 },
    @{ Cell = "F22"; Value = 0.9916831474676431 },
    @{ Cell = "J22"; Value = 1.003999810658025 },
    @{ Cell = "K22"; Value = 1.004416738918358 },
    @{ Cell = "L22"; Value = 1.00467084750642 },
    @{ Cell = "M22"; Value = 0.9955720865436291 },
    @{ Cell = "N22"; Value = 1.004938784099908 },
    @{ Cell = "C23"; Value = 0.9976617515572034 },
    @{ Cell = "D23"; Value = 1.001285068720394 },
    @{ Cell = "E23"; Value = 1.001513637737119 },
    @{ Cell = "F23"; Value = 0.9927401879471073 },
    @{ Cell = "J23"; Value = 1.004656384205627 },
    @{ Cell = "K23"; Value = 1.005062507507651 },
    @{ Cell = "L23"; Value = 1.005290142572133 },
    @{ Cell = "M23"; Value = 0.9965528862430468 },
    @{ Cell = "N23"; Value = 1.00517365239791 },
    @{ Cell = "C24"; Value = 1.000789611011587 },
    @{ Cell = "D24"; Value = 1.004104493814722 },
    @{ Cell = "E24"; Value = 1.004228773359167 },
    @{ Cell = "F24"; Value = 0.9968812494534708 },
    @{ Cell = "J24"; Value = 1.007225743100768 },
    @{ Cell = "K24"; Value = 1.007589768787543 },
    @{ Cell = "L24"; Value = 1.007713583694287 },
    @{ Cell = "M24"; Value = 1.000393750163362 },
    @{ Cell = "N24"; Value = 1.006090756012724 },
    @{ Cell = "C25"; Value = 1.004387773426811 },
    @{ Cell = "D25"; Value = 1.007347887935683 },
    @{ Cell = "E25"; Value = 1.007352013251212 },
    @{ Cell = "F25"; Value = 1.001646885398629 },
    @{ Cell = "J25"; Value = 1.010176558077595 },
    @{ Cell = "K25"; Value = 1.010492603475837 },
    @{ Cell = "L25"; Value = 1.010496715010997 },
    @{ Cell = "M25"; Value = 1.004810763952434 },
    @{ Cell = "N25"; Value = 1.007139822571723 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value2 = $u.Value
}
